# "Elimina EC anteriores y se agregan nuevos, se modifica base de datos"
#
# The "Periodo Mora" list (column E, rows 16-22) is refreshed: the periods
# are re-entered in descending order (2209 down to 2203) instead of the
# previous ascending order (2203 up to 2209), and each period's
# "Valor Mora" (column F) travels with its period so the pairing between a
# period and its mora value is preserved. Columns B, C, D, G, H, I, J are
# identical across this block and are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 16
$lastRow = 22
$n = $lastRow - $firstRow + 1

# Read the current Periodo Mora (col E) / Valor Mora (col F) pairs.
# NOTE: use .Value2 for reads -- .Value as a getter is unreliable in this
# host and .Value2 round-trips numbers/strings cleanly.
$periods = @()
$values = @()
for ($i = 0; $i -lt $n; $i++) {
    $r = $firstRow + $i
    $periods += $ws.Cells.Item($r, 5).Value2
    $values += $ws.Cells.Item($r, 6).Value2
}

# Sort the (period, value) pairs by period, descending -- a simple
# in-place bubble sort since Sort-Object isn't available in this host.
for ($i = 0; $i -lt $n; $i++) {
    for ($j = 0; $j -lt ($n - $i - 1); $j++) {
        if ($periods[$j] -lt $periods[$j + 1]) {
            $tmpP = $periods[$j]
            $periods[$j] = $periods[$j + 1]
            $periods[$j + 1] = $tmpP

            $tmpV = $values[$j]
            $values[$j] = $values[$j + 1]
            $values[$j + 1] = $tmpV
        }
    }
}

# Write the reordered pairs back into the same rows.
for ($i = 0; $i -lt $n; $i++) {
    $r = $firstRow + $i
    $ws.Cells.Item($r, 5).Value = $periods[$i]
    $ws.Cells.Item($r, 6).Value = $values[$i]
}
